$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.146.00'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '''1.673.25'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("D5").Value = '''216.98'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '''0.5220'
$ws.Range("E6").Value = '  +1.81%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = '''0.2701'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").Value = '''0.06395'
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("D11").Value = '''0.07434'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '''1.688.47'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").Value = '''4.521'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '''0.5833'
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '''0.000008532'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '''64.26'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '''25.944.39'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").Value = '''4.933'
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '''10.80'
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = '''189.92'
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").Value = '''6.191'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = '''1.004'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '''144.89'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").Value = '''0.1243'
$ws.Range("E25").Value = '  +5.40%  '
$ws.Range("D26").Value = '''7.618'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '''15.71'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '''0.06574'
$ws.Range("E28").Value = '  +13.16%  '
$ws.Range("D29").Value = '''1.337'
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").Value = '''1.319'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("D32").Value = '''3.533'
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").Value = '''1.668'
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("D35").Value = '''0.6166'
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''2.700'
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").Value = '''6.278'
$ws.Range("E38").Value = '  +6.38%  '
$ws.Range("D39").Value = '''1.094.76'
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").Value = '''0.01598'
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").Value = '''0.8718'
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").Value = '''100.83'
$ws.Range("D44").Value = '''1.819.30'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '''0.00000000110'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("D46").Value = '''56.51'
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("D47").Value = '''8.159'
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("D48").Value = '''1.001'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").Value = '''0.05240'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '''0.4280'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").Value = '''5.995'
$ws.Range("E51").Value = '  +2.58%  '
